$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "HongKong"
$ws.Range("B30").Value = "Issue with OSM rules"

$ws.Range("A31").Value = "Azerbeijan "
$ws.Range("B31").Value = "Issues with OSM cables"

$ws.Range("A32").Value = "Laos"
$ws.Range("B32").Value = "Issues with OSM cables"

$ws.Range("A33").Value = "Puerto Rico"
$ws.Range("B33").Value = "No powerplants so ppmatching breaks"

$ws.Range("A34").Value = "New Zealand"
$ws.Range("B34").Value = "Cutout seems to be too large. Need to look into this"

$ws.Range("A35").Value = "Sri Lanka"
$ws.Range("B35").Value = "Issues with OSM cables"

$ws.Range("A36").Value = "Panama"
$ws.Range("B36").Value = "Issues with OSM cables"

$ws.Range("A37").Value = "Armenia"
$ws.Range("B37").Value = "Issues with OSM cables"

$ws.Range("A30:B37").WrapText = $true
$ws.Rows("30:37").AutoFit()

$ws.Range("A30").Select()
